$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 725.25
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 725.25
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 725.25
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1377.25

$ws.Range("H112").Value = 6201.7334
$ws.Range("I112").Value = 35250
$ws.Range("J112").Value = 1732.7693
$ws.Range("K112").Value = 105750
$ws.Range("L112").Value = 5198.3079
$ws.Range("M112").Value = -104642
$ws.Range("N112").Value = -7414.3079

$ws.Range("H113").Value = 3504.8
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3881
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3881
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -10389

$ws.Range("H128").Value = 37475
$ws.Range("J128").Value = 37475
$ws.Range("L128").Value = 37475
$ws.Range("N128").Value = -47435

$ws.Range("H138").Value = 9924.489
$ws.Range("I138").Value = 2681.6365
$ws.Range("J138").Value = 12267.765
$ws.Range("K138").Value = 8044.9095
$ws.Range("L138").Value = 36803.295
$ws.Range("M138").Value = -2904.9095
$ws.Range("N138").Value = -47083.295

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 5000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -5458

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H27").Value = 34000
$ws.Range("J27").Value = 34000
$ws.Range("L27").Value = 34000
$ws.Range("N27").Value = -34368

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H32").Value = 2581.55
$ws.Range("I32").Value = 2639.1353
$ws.Range("K32").Value = 2639.1353
$ws.Range("M32").Value = -2352.1353

$ws.Range("H74").Value = 6150.8076
$ws.Range("I74").Value = 2325.2632
$ws.Range("J74").Value = 16534.428
$ws.Range("K74").Value = 2325.2632
$ws.Range("L74").Value = 16534.428
$ws.Range("M74").Value = -1451.2632
$ws.Range("N74").Value = -18282.428

$ws.Range("H77").Value = 6150.8076
$ws.Range("I77").Value = 2325.2632
$ws.Range("J77").Value = 16534.428
$ws.Range("K77").Value = 11626.316
$ws.Range("L77").Value = 82672.14
$ws.Range("M77").Value = -7258.315999999999
$ws.Range("N77").Value = -91408.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 64714.285
$ws.Range("J117").Value = 64714.285
$ws.Range("L117").Value = 64714.285
$ws.Range("N117").Value = -73892.285

$ws.Range("H132").Value = 9254.823
$ws.Range("I132").Value = 11551.667
$ws.Range("J132").Value = 3742.4
$ws.Range("K132").Value = 34655.001
$ws.Range("L132").Value = 11227.2
$ws.Range("M132").Value = -32125.001
$ws.Range("N132").Value = -16287.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3142.375
$ws.Range("I70").Value = 700
$ws.Range("J70").Value = 3305.2
$ws.Range("K70").Value = 2100
$ws.Range("L70").Value = 9915.599999999999
$ws.Range("M70").Value = -1785
$ws.Range("N70").Value = -10545.6

$ws.Range("H73").Value = 3142.375
$ws.Range("I73").Value = 700
$ws.Range("J73").Value = 3305.2
$ws.Range("K73").Value = 2100
$ws.Range("L73").Value = 9915.599999999999
$ws.Range("M73").Value = -1008
$ws.Range("N73").Value = -12099.6

$ws.Range("H129").Value = 1724.4615
$ws.Range("I129").Value = 1551.5385
$ws.Range("J129").Value = 1897.3846
$ws.Range("K129").Value = 4654.6155
$ws.Range("L129").Value = 5692.1538
$ws.Range("M129").Value = 345.3845000000001
$ws.Range("N129").Value = -15692.1538

$ws.Range("H137").Value = 60077.332
$ws.Range("J137").Value = 107173.2
$ws.Range("L137").Value = 321519.6
$ws.Range("N137").Value = -331719.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 3401.3333
$ws.Range("I9").Value = 2300
$ws.Range("J9").Value = 4502.6665
$ws.Range("K9").Value = 2300
$ws.Range("L9").Value = 4502.6665
$ws.Range("M9").Value = -2130
$ws.Range("N9").Value = -4842.6665

$ws.Range("H24").Value = 71714.28999999999
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15346

$ws.Range("H62").Value = 41723.332
$ws.Range("J62").Value = 41723.332
$ws.Range("L62").Value = 41723.332
$ws.Range("N62").Value = -43095.332

$ws.Range("H65").Value = 41723.332
$ws.Range("J65").Value = 41723.332
$ws.Range("L65").Value = 125169.996
$ws.Range("N65").Value = -132033.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3649.5881
$ws.Range("I7").Value = 3377.375
$ws.Range("K7").Value = 3377.375
$ws.Range("M7").Value = -3265.375

$ws.Range("H21").Value = 25018750
$ws.Range("I21").Value = 15000
$ws.Range("J21").Value = 33353334
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 33353334
$ws.Range("M21").Value = -14826
$ws.Range("N21").Value = -33353682

$ws.Range("H22").Value = 590
$ws.Range("I22").Value = 590
$ws.Range("K22").Value = 590
$ws.Range("M22").Value = -295

$ws.Range("H27").Value = 590
$ws.Range("I27").Value = 590
$ws.Range("K27").Value = 590
$ws.Range("M27").Value = -483

$ws.Range("H46").Value = 1166.5
$ws.Range("I46").Value = 1533
$ws.Range("J46").Value = 800
$ws.Range("K46").Value = 1533
$ws.Range("L46").Value = 800
$ws.Range("M46").Value = -1345
$ws.Range("N46").Value = -1176

$ws.Range("H68").Value = 4450
$ws.Range("J68").Value = 4300
$ws.Range("L68").Value = 4300
$ws.Range("N68").Value = -5798

$ws.Range("H71").Value = 4450
$ws.Range("J71").Value = 4300
$ws.Range("L71").Value = 21500
$ws.Range("N71").Value = -28988

$ws.Range("H126").Value = 3649.5881
$ws.Range("I126").Value = 3377.375
$ws.Range("K126").Value = 10132.125
$ws.Range("M126").Value = -7662.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 7933.3335
$ws.Range("J29").Value = 7933.3335
$ws.Range("L29").Value = 7933.3335
$ws.Range("N29").Value = -8513.333500000001

$ws.Range("H109").Value = 66266.664
$ws.Range("J109").Value = 66266.664
$ws.Range("L109").Value = 66266.664
$ws.Range("N109").Value = -69040.664

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 3741.5417
$ws.Range("I132").Value = 3325.875
$ws.Range("J132").Value = 4572.875
$ws.Range("K132").Value = 9977.625
$ws.Range("L132").Value = 13718.625
$ws.Range("M132").Value = -7447.625
$ws.Range("N132").Value = -18778.625
